$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from
# 45175 (2023-09-06) to 45177 (2023-09-08) for every data row (2-99).
for ($row = 2; $row -le 99; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
